$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated cryptocurrency price (D) and volume/1h-change (E) figures.
# D-column values are written through a text-format/style round trip so that
# numeric-looking strings (e.g. "9.60", "0.999") are preserved verbatim as text
# (matching the original inline-string cell type) instead of being auto-coerced
# into numbers by Excel (which would silently drop significant trailing zeros).

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = '71.385.59'
$c.Style = "Normal"
$ws.Range("E2").Value = '  +3.35%  '

$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = '3.584.31'
$c.Style = "Normal"
$ws.Range("E3").Value = '  +1.84%  '

$ws.Range("E4").Value = '  -0.04%  '

$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = '584.64'
$c.Style = "Normal"
$ws.Range("E5").Value = '  +2.49%  '

$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = '187.72'
$c.Style = "Normal"
$ws.Range("E6").Value = '  +3.50%  '

$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = '3.574.70'
$c.Style = "Normal"
$ws.Range("E7").Value = '  +1.71%  '

$ws.Range("E8").Value = '  +1.81%  '

$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = '0.999'
$c.Style = "Normal"
$ws.Range("E9").Value = '  -0.15%  '

$ws.Range("E10").Value = '  +14.76%  '

$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = '0.656'
$c.Style = "Normal"
$ws.Range("E11").Value = '  +3.24%  '

$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = '54.81'
$c.Style = "Normal"
$ws.Range("E12").Value = '  +2.63%  '

$ws.Range("E13").Value = '  +5.94%  '

$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = '9.60'
$c.Style = "Normal"
$ws.Range("E14").Value = '  +1.93%  '

$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = '4.056.23'
$c.Style = "Normal"
$ws.Range("E15").Value = '  -1.14%  '

$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = '71.322.12'
$c.Style = "Normal"
$ws.Range("E16").Value = '  +3.27%  '

$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = '19.33'
$c.Style = "Normal"
$ws.Range("E17").Value = '  +0.67%  '

$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = '3.581.09'
$c.Style = "Normal"
$ws.Range("E18").Value = '  +1.75%  '

$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = '12.42'
$c.Style = "Normal"
$ws.Range("E19").Value = '  +0.12%  '

$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = '566.82'
$c.Style = "Normal"
$ws.Range("E20").Value = '  +5.56%  '

$ws.Range("E21").Value = '  +0.71%  '

$ws.Range("E22").Value = '  -1.12%  '

$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = '17.59'
$c.Style = "Normal"
$ws.Range("E23").Value = '  -11.04%  '

$ws.Range("E24").Value = '  +2.54%  '

$ws.Range("E25").Value = '  +5.24%  '

$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = '94.89'
$c.Style = "Normal"
$ws.Range("E26").Value = '  +0.83%  '

$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = '11.34'
$c.Style = "Normal"
$ws.Range("E27").Value = '  +3.27%  '

$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$ws.Range("E28").Value = '  +2.28%  '

$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = '9.22'
$c.Style = "Normal"
$ws.Range("E29").Value = '  +2.08%  '

$ws.Range("E30").Value = '  +3.51%  '

$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = '7.33'
$c.Style = "Normal"
$ws.Range("E31").Value = '  +1.24%  '

$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = '12.40'
$c.Style = "Normal"
$ws.Range("E32").Value = '  -0.72%  '

$ws.Range("E33").Value = '  +1.78%  '

$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = '64.22'
$c.Style = "Normal"
$ws.Range("E34").Value = '  -0.81%  '

$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = '3.43'
$c.Style = "Normal"
$ws.Range("E35").Value = '  +9.33%  '

$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = '555.41'
$c.Style = "Normal"
$ws.Range("E36").Value = '  -2.28%  '

$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = '0.423'
$c.Style = "Normal"
$ws.Range("E37").Value = '  +6.79%  '

$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = '0.0₃0807'
$c.Style = "Normal"
$ws.Range("E38").Value = '  +6.22%  '

$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = '37.85'
$c.Style = "Normal"
$ws.Range("E39").Value = '  -0.22%  '

$ws.Range("E40").Value = '  +0.20%  '

$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = '3.28'
$c.Style = "Normal"
$ws.Range("E41").Value = '  +9.29%  '

$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = '3.527.93'
$c.Style = "Normal"
$ws.Range("E42").Value = '  +11.76%  '

$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = '3.46'
$c.Style = "Normal"
$ws.Range("E43").Value = '  +3.19%  '

$ws.Range("E44").Value = '  +3.43%  '

$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = '0.0451'
$c.Style = "Normal"
$ws.Range("E45").Value = '  +2.97%  '

$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = '2.96'
$c.Style = "Normal"
$ws.Range("E46").Value = '  +0.61%  '

$ws.Range("E47").Value = '  -0.99%  '

$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = '9.41'
$c.Style = "Normal"
$ws.Range("E48").Value = '  +2.30%  '

$ws.Range("E49").Value = '  +3.56%  '

$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = '1.49'
$c.Style = "Normal"
$ws.Range("E50").Value = '  +8.37%  '

$ws.Range("E51").Value = '  -0.10%  '
